# Fill in previously-missing stock data (Closing Price / Value) for the
# Motilal Oswal MOSt Shares NASDAQ-100 ETF (MON100) rows, which used to be
# left blank (fetched from server on demand). Now that this data is always
# present in the prefetched data, these cells get real numeric values, and
# every Proportion (H) cell in the sheet is recomputed to account for the
# newly-included MON100 values in each day's total portfolio value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H15").Value = 0.1897198552766218

$ws.Range("H16").Value = 0.516962756928909

$ws.Range("F17").Value = 615.3300094604492
$ws.Range("G17").Value = 3076.650047302246
$ws.Range("H17").Value = 0.2933173877944692

$ws.Range("H18").Value = 0.1963294899032366

$ws.Range("H19").Value = 0.471190775767768

$ws.Range("F20").Value = 616.5399932861328
$ws.Range("G20").Value = 3082.699966430664
$ws.Range("H20").Value = 0.2392193399940886

$ws.Range("H21").Value = 0.09326039433490677

$ws.Range("H22").Value = 0.2342788964336497

$ws.Range("H23").Value = 0.4123393302057723

$ws.Range("F24").Value = 628.1700134277344
$ws.Range("G24").Value = 3769.020080566406
$ws.Range("H24").Value = 0.2655245458626109

$ws.Range("H25").Value = 0.08785722749796712

$ws.Range("H26").Value = 0.2200738230840038

$ws.Range("H27").Value = 0.3831605434493547

$ws.Range("F28").Value = 641.4800262451172
$ws.Range("G28").Value = 4490.36018371582
$ws.Range("H28").Value = 0.3101728691311063

$ws.Range("H29").Value = 0.0865927643355353

$ws.Range("H30").Value = 0.2105908844293644

$ws.Range("H31").Value = 0.343445404846221

$ws.Range("F32").Value = 637.3600006103516
$ws.Range("G32").Value = 6373.600006103516
$ws.Range("H32").Value = 0.3732282411635292

$ws.Range("H33").Value = 0.0727354695608855

$ws.Range("H34").Value = 0.3066465430776691

$ws.Range("F35").Value = 649.27099609375
$ws.Range("G35").Value = 6492.7099609375
$ws.Range("H35").Value = 0.5792451775289043

$ws.Range("H36").Value = 0.1141082793934266

$ws.Range("H37").Value = 0.2756603970605329

$ws.Range("F38").Value = 666.5799713134766
$ws.Range("G38").Value = 6665.799713134766
$ws.Range("H38").Value = 0.5537965450563179

$ws.Range("H39").Value = 0.1134045900767662

$ws.Range("H40").Value = 0.05713846780638292

$ws.Range("H41").Value = 0.2675888384137958

$ws.Range("F42").Value = 662.6899719238281
$ws.Range("G42").Value = 6626.899719238281
$ws.Range("H42").Value = 0.4764331918316316

$ws.Range("H43").Value = 0.09540311040548617

$ws.Range("H44").Value = 0.1605748593490865

$ws.Range("H45").Value = 0.1957773201390204

$ws.Range("F46").Value = 682.9499816894531
$ws.Range("G46").Value = 10244.2497253418
$ws.Range("H46").Value = 0.5773148512662105

$ws.Range("H47").Value = 0.07288957832364869

$ws.Range("H48").Value = 0.1540182502711204

$ws.Range("H49").Value = 0.1605961443530085

$ws.Range("F50").Value = 630.0400161743164
$ws.Range("G50").Value = 10710.68027496338
$ws.Range("H50").Value = 0.5699450050133548

$ws.Range("H51").Value = 0.06105101779666405

$ws.Range("H52").Value = 0.2084078328369726

$ws.Range("H53").Value = 0.1182164291121852

$ws.Range("F54").Value = 526.7499923706055
$ws.Range("G54").Value = 10534.99984741211
$ws.Range("H54").Value = 0.6798089974915936

$ws.Range("H55").Value = 0.05301026098984984

$ws.Range("H56").Value = 0.1489643124063713
